$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 89; this shifts the existing rows 89..187
# down to 91..189, preserving all their data and formatting.
$ws.Rows("89:90").Insert()

# ---- New row 89: Clementina / Primera (week of 44790) ----
$ws.Cells.Item(89,1).Value  = 7
$ws.Cells.Item(89,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(89,3).Value  = "Ñuble"
$ws.Cells.Item(89,4).Value  = 44790
$ws.Cells.Item(89,5).Value  = 16
$ws.Cells.Item(89,6).Value  = "Fruta"
$ws.Cells.Item(89,7).Value  = 100102
$ws.Cells.Item(89,8).Value  = "Cítricos"
$ws.Cells.Item(89,9).Value  = 100102004
$ws.Cells.Item(89,10).Value = "Mandarina"
$ws.Cells.Item(89,11).Value = "Clementina"
$ws.Cells.Item(89,12).Value = "Primera"
$ws.Cells.Item(89,13).Value = 80
$ws.Cells.Item(89,14).Value = 8500
$ws.Cells.Item(89,15).Value = 9000
$ws.Cells.Item(89,16).Value = 8750
$ws.Cells.Item(89,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(89,18).Value = "Región de O'Higgins"
$ws.Cells.Item(89,19).Value = 486
$ws.Cells.Item(89,20).Value = 18

# ---- New row 90: Clementina / Segunda (week of 44790) ----
$ws.Cells.Item(90,1).Value  = 7
$ws.Cells.Item(90,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(90,3).Value  = "Ñuble"
$ws.Cells.Item(90,4).Value  = 44790
$ws.Cells.Item(90,5).Value  = 16
$ws.Cells.Item(90,6).Value  = "Fruta"
$ws.Cells.Item(90,7).Value  = 100102
$ws.Cells.Item(90,8).Value  = "Cítricos"
$ws.Cells.Item(90,9).Value  = 100102004
$ws.Cells.Item(90,10).Value = "Mandarina"
$ws.Cells.Item(90,11).Value = "Clementina"
$ws.Cells.Item(90,12).Value = "Segunda"
$ws.Cells.Item(90,13).Value = 120
$ws.Cells.Item(90,14).Value = 7500
$ws.Cells.Item(90,15).Value = 8000
$ws.Cells.Item(90,16).Value = 7750
$ws.Cells.Item(90,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(90,18).Value = "Región de O'Higgins"
$ws.Cells.Item(90,19).Value = 431
$ws.Cells.Item(90,20).Value = 18
